$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("R30" rule) column C ("From") value is corrected from 18 to 1.
$ws.Range("C10").Value = 1
